# 1. add pi and theta into optimization model 2. income process simulation
#
# Applies:
#  - workbook window position (best effort; host window geometry)
#  - "Coefficients" sheet: replace formula-derived coefficients with plain
#    simulated values (pi/theta optimization inputs), update selection
#  - "Variance" sheet: highlight (yellow fill) the sigma-u/sigma-e point
#    estimates for the "Labor Income Only" columns, update selection

$wb = $excel.ActiveWorkbook

# --- Workbook window position (best effort) ---
$win = $excel.ActiveWindow
$win.Left = 1160
$win.Top = 1920

# --- Sheet: Coefficients ---
$wsCoef = $wb.Worksheets.Item("Coefficients")

$wsCoef.Range("B2").Value = 41078.88
$wsCoef.Range("C2").Value = -2201.582
$wsCoef.Range("D2").Value = 78.22341
$wsCoef.Range("E2").Value = -0.7889483

$wsCoef.Range("B3").Value = 28752.76
$wsCoef.Range("C3").Value = -917.5466
$wsCoef.Range("D3").Value = 51.09514
$wsCoef.Range("E3").Value = -0.6089281

[void]$wsCoef.Range("C11").Select()

# --- Sheet: Variance ---
$wsVar = $wb.Worksheets.Item("Variance")
$wsVar.Activate()

$wsVar.Range("B4").Interior.Color = 65535
$wsVar.Range("C4").Interior.Color = 65535
$wsVar.Range("B6").Interior.Color = 65535
$wsVar.Range("C6").Interior.Color = 65535

[void]$wsVar.Range("E2:G2").Select()
